$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Check_Points": new row 12 entry -- a hyperlinked "Benefits Online"
# link plus its URL shown as plain text in column C.
# ---------------------------------------------------------------------------
$wsCheck = $wb.Worksheets.Item("Check_Points")

# Start from a clean (unstyled) cell so the new cellXfs entry Excel creates
# for the hyperlink only carries the Hyperlink font, not whatever
# border/alignment happened to live on the cell beforehand.
$wsCheck.Range("B12").ClearFormats()
$wsCheck.Range("B12").Value = "Benefits Online | Sell Shares"

$wsCheck.Hyperlinks.Add($wsCheck.Range("B12"), "https://www.benefits.ml.com/Core/Frame/ContentHostV2?MenuID=11066", "", "https://www.benefits.ml.com/Core/Frame/ContentHostV2?MenuID=11066", "https://www.benefits.ml.com/Core/Frame/ContentHostV2?MenuID=11066")

# Restore the friendly display text in B12 (kept as its own shared string)
# and put the raw URL as plain text in C12.
$wsCheck.Range("B12").Value = "Benefits Online | Sell Shares"
$wsCheck.Range("C12").Value = "https://www.benefits.ml.com/Core/Frame/ContentHostV2?MenuID=11066"

# ---------------------------------------------------------------------------
# Sheet "Tracking": the Date/Day columns (B:C) shift up five rows -- the
# first five days' worth of entries were removed, so row 2 now shows what
# used to be row 7, and so on, leaving the last five rows of B:C blank.
# ---------------------------------------------------------------------------
$wsTrack = $wb.Worksheets.Item("Tracking")

for ($r = 2; $r -le 19; $r++) {
    $src = $r + 5
    $wsTrack.Cells.Item($r, 2).Value2 = $wsTrack.Cells.Item($src, 2).Value2
    $wsTrack.Cells.Item($r, 3).Value2 = $wsTrack.Cells.Item($src, 3).Value2
}
for ($r = 20; $r -le 24; $r++) {
    $wsTrack.Cells.Item($r, 2).Clear()
    $wsTrack.Cells.Item($r, 3).Clear()
}

# ---------------------------------------------------------------------------
# View state: Check_Points becomes the active/selected tab (cursor on A12),
# Tracking is no longer the active tab (cursor moves to E12), Topics is
# left untouched.
# ---------------------------------------------------------------------------
$wsTrack.Activate()
$wsTrack.Range("E12").Select()

$wsCheck.Activate()
$wsCheck.Range("A12").Select()
